$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New conversation rows appended to the media-conversation log (rows 43-55).
# Columns: A=Timestamp, B=Sender, C=Sender Id (numeric), D=Phone (text),
#          E=Message, F=Media (text, blank for some rows), G=Channel (blank)

$rows = @(
    @{ Row=43; Timestamp="2025-10-02 20:29:18"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="This is a test message"; Media="my-node-server/public/uploads/images\photo_2025-10-03_00-29-18.jpg" }
    @{ Row=44; Timestamp="2025-10-02 20:34:00"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="This is a test message"; Media="my-node-server/public/uploads/images\photo_2025-10-03_00-34-01.jpg" }
    @{ Row=45; Timestamp="2025-10-02 20:34:46"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="This is a test message"; Media="" }
    @{ Row=46; Timestamp="2025-10-02 20:38:35"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="This is a test message"; Media="" }
    @{ Row=47; Timestamp="2025-10-02 20:40:25"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="What’s up?"; Media="" }
    @{ Row=48; Timestamp="2025-10-02 20:40:44"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="This is a test message"; Media="my-node-server/public/uploads/images\photo_2025-10-03_00-40-45.jpg" }
    @{ Row=49; Timestamp="2025-10-02 20:51:16"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="Test message"; Media="" }
    @{ Row=50; Timestamp="2025-10-02 21:09:20"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="This is a test message"; Media="my-node-server/public/uploads/images\photo_2025-10-03_01-09-20.jpg" }
    @{ Row=51; Timestamp="2025-10-02 21:10:27"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="Test message"; Media="" }
    @{ Row=52; Timestamp="2025-10-02 21:11:06"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="Latest message"; Media="" }
    @{ Row=53; Timestamp="2025-10-02 21:11:17"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="Message again"; Media="" }
    @{ Row=54; Timestamp="2025-10-02 21:12:37"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="Test message"; Media="" }
    @{ Row=55; Timestamp="2025-10-02 21:13:09"; Sender="Noah"; SenderId=8450689526; Phone="13052054965"; Message="Test 113"; Media="" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Timestamp
    $ws.Cells.Item($row, 2).Value = $r.Sender
    $ws.Cells.Item($row, 3).Value = $r.SenderId
    # Leading apostrophe forces the numeric-looking phone number to be
    # stored as text, matching the rest of column D in the sheet.
    $ws.Cells.Item($row, 4).Value = "'" + $r.Phone
    $ws.Cells.Item($row, 5).Value = $r.Message
    if ($r.Media -ne "") {
        $ws.Cells.Item($row, 6).Value = $r.Media
    }
}
